$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44, 188),
    @(47, 198),
    @(40, 184),
    @(42, 173),
    @(45, 180),
    @(40, 169)
)

$r = 7
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r++
}

$ws.Range("D8").Select()
